$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column F, copying the header style (bold/border/centered) from E1
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the time_taken values for each data row (2-49)
$timeTaken = @(
    "2021-10-05 10:51:32.728590",
    "2021-10-05 10:51:32.728603",
    "2021-10-05 10:51:32.728607",
    "2021-10-05 10:51:32.728610",
    "2021-10-05 10:51:32.728614",
    "2021-10-05 10:51:32.728617",
    "2021-10-05 10:51:32.728620",
    "2021-10-05 10:51:32.728623",
    "2021-10-05 10:51:32.728626",
    "2021-10-05 10:51:32.728630",
    "2021-10-05 10:51:32.728633",
    "2021-10-05 10:51:32.728635",
    "2021-10-05 10:51:32.728638",
    "2021-10-05 10:51:32.728641",
    "2021-10-05 10:51:32.728644",
    "2021-10-05 10:51:32.728647",
    "2021-10-05 10:51:32.728651",
    "2021-10-05 10:51:32.728654",
    "2021-10-05 10:51:32.728657",
    "2021-10-05 10:51:32.728660",
    "2021-10-05 10:51:32.728663",
    "2021-10-05 10:51:32.728666",
    "2021-10-05 10:51:32.728669",
    "2021-10-05 10:51:32.728672",
    "2021-10-05 10:51:32.728675",
    "2021-10-05 10:51:32.728678",
    "2021-10-05 10:51:32.728681",
    "2021-10-05 10:51:32.728684",
    "2021-10-05 10:51:32.728687",
    "2021-10-05 10:51:32.728690",
    "2021-10-05 10:51:32.728693",
    "2021-10-05 10:51:32.728696",
    "2021-10-05 10:51:32.728700",
    "2021-10-05 10:51:32.728703",
    "2021-10-05 10:51:32.728706",
    "2021-10-05 10:51:32.728709",
    "2021-10-05 10:51:32.728712",
    "2021-10-05 10:51:32.728715",
    "2021-10-05 10:51:32.728719",
    "2021-10-05 10:51:32.728721",
    "2021-10-05 10:51:32.728725",
    "2021-10-05 10:51:32.728728",
    "2021-10-05 10:51:32.728731",
    "2021-10-05 10:51:32.728734",
    "2021-10-05 10:51:32.728737",
    "2021-10-05 10:51:32.728740",
    "2021-10-05 10:51:32.728743",
    "2021-10-05 10:51:32.728746"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}

